# finished link skills, added MATT
# Inserts a new "matt" (magic attack) column before the existing "def" column (L),
# shifting def/bossdamage/ied/finaldamage/flathp/flatmp/%hp/%mp one column right,
# then fills in the new column's header and per-row magic-attack values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at L (shifts L:S -> M:T)
$ws.Range("L1").EntireColumn.Insert()

# New header
$ws.Range("L1").Value = "matt"

# New per-row magic attack values (rows without att/def data are left blank)
$mattValues = @{
  2  = 0
  3  = 3
  4  = 10
  5  = 77
  6  = 3
  7  = 10
  8  = 10
  9  = 10
  10 = 10
  11 = 5
  12 = 4
  13 = 6
  14 = 10
  15 = 5
  16 = 5
  17 = 5
  18 = 5
  19 = 5
  20 = 0
  21 = 0
  22 = 0
  23 = 0
  24 = 0
  25 = 0
  26 = 0
  27 = 0
  28 = 0
  29 = 0
  30 = 0
  31 = 0
  32 = 0
  33 = 0
  34 = 0
  35 = 0
  36 = 0
  37 = 0
  38 = 0
  39 = 0
  40 = 0
}

foreach ($row in $mattValues.Keys) {
    $ws.Range("L$row").Value = $mattValues[$row]
}

# Match the author's final selection from the commit
$ws.Range("L54").Select()
